$d = $word.ActiveDocument

function Get-ParaIndexByText($needle) {
    $r = $d.Content
    $found = $r.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $needle"
    }
    $pos = $r.Start
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($pos -ge $p.Range.Start -and $pos -lt $p.Range.End) {
            return $i
        }
    }
    throw "Could not resolve paragraph index for: $needle"
}

# ------------------------------------------------------------------
# 1. Delete the original "3.1.4 日志配置" section. It currently sits
#    right after the 3.1.3 code-generation section (after "...生成的
#    表要有注释"), and right before the "注：如对模板有特殊需求..."
#    note. Remove it together with its leading/trailing blank spacer
#    paragraphs.
# ------------------------------------------------------------------
$headingIdx = Get-ParaIndexByText("日志配置")
$startPara = $d.Paragraphs.Item($headingIdx).Previous(1)
$endPara = $d.Paragraphs.Item($headingIdx)
while ($true) {
    $nextPara = $endPara.Next(1)
    $t = $nextPara.Range.Text
    if ($t -match "注：如对模板") { break }
    $endPara = $nextPara
}

$delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$delRange.Delete()

# ------------------------------------------------------------------
# 2. Re-insert the (reformatted) section right after the screenshot
#    that follows the "templates/vm" note, i.e. directly before the
#    "3.2 启动及验证" heading.
# ------------------------------------------------------------------
$targetIdx = Get-ParaIndexByText("3.2")
$targetPara = $d.Paragraphs.Item($targetIdx)
$insertionRange = $targetPara.Range
$insertionRange.Collapse(1)
$insertionRange.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item($targetIdx).Previous(1)
$newRange = $newPara.Range

$xmlPayload = @'
<w:p><w:pPr><w:pStyle w:val="10"/><w:jc w:val="left"/><w:rPr><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve">3.1.4 </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>日志配置</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:szCs w:val="21"/></w:rPr><w:t>编辑</w:t></w:r><w:r><w:rPr><w:szCs w:val="21"/></w:rPr><w:t>src/</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:szCs w:val="21"/></w:rPr><w:t>main</w:t></w:r><w:r><w:rPr><w:szCs w:val="21"/></w:rPr><w:t>/</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:szCs w:val="21"/></w:rPr><w:t>resources</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:szCs w:val="21"/></w:rPr><w:t>目录下的</w:t></w:r><w:r><w:rPr><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve">logback.yml </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:szCs w:val="21"/></w:rPr><w:t>文件</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="11"/><w:keepNext w:val="0"/><w:keepLines w:val="0"/><w:widowControl/><w:suppressLineNumbers w:val="0"/><w:shd w:val="clear" w:fill="FFFFFF"/><w:rPr><w:rFonts w:ascii="Menlo" w:hAnsi="Menlo" w:eastAsia="Menlo" w:cs="Menlo"/><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Menlo" w:hAnsi="Menlo" w:eastAsia="Menlo" w:cs="Menlo"/><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:shd w:val="clear" w:fill="EFEFEF"/></w:rPr><w:t>&lt;</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Menlo" w:hAnsi="Menlo" w:eastAsia="Menlo" w:cs="Menlo"/><w:b/><w:color w:val="000080"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:shd w:val="clear" w:fill="EFEFEF"/></w:rPr><w:t xml:space="preserve">property </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Menlo" w:hAnsi="Menlo" w:eastAsia="Menlo" w:cs="Menlo"/><w:b/><w:color w:val="0000FF"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:shd w:val="clear" w:fill="EFEFEF"/></w:rPr><w:t>name</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Menlo" w:hAnsi="Menlo" w:eastAsia="Menlo" w:cs="Menlo"/><w:b/><w:color w:val="008000"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:shd w:val="clear" w:fill="EFEFEF"/></w:rPr><w:t xml:space="preserve">="log.path" </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Menlo" w:hAnsi="Menlo" w:eastAsia="Menlo" w:cs="Menlo"/><w:b/><w:color w:val="0000FF"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:shd w:val="clear" w:fill="EFEFEF"/></w:rPr><w:t>value</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Menlo" w:hAnsi="Menlo" w:eastAsia="Menlo" w:cs="Menlo"/><w:b/><w:color w:val="008000"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:shd w:val="clear" w:fill="EFEFEF"/></w:rPr><w:t xml:space="preserve">="/home/ruoyi/logs" </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Menlo" w:hAnsi="Menlo" w:eastAsia="Menlo" w:cs="Menlo"/><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:shd w:val="clear" w:fill="EFEFEF"/></w:rPr><w:t>/&gt;</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default"/><w:szCs w:val="21"/></w:rPr><w:t>改为自己需要的路径</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="11"/><w:keepNext w:val="0"/><w:keepLines w:val="0"/><w:widowControl/><w:suppressLineNumbers w:val="0"/><w:shd w:val="clear" w:fill="FFFFFF"/><w:rPr><w:rFonts w:ascii="Menlo" w:hAnsi="Menlo" w:eastAsia="Menlo" w:cs="Menlo"/><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Menlo" w:hAnsi="Menlo" w:eastAsia="Menlo" w:cs="Menlo"/><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:shd w:val="clear" w:fill="EFEFEF"/></w:rPr><w:t>&lt;</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Menlo" w:hAnsi="Menlo" w:eastAsia="Menlo" w:cs="Menlo"/><w:b/><w:color w:val="000080"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:shd w:val="clear" w:fill="EFEFEF"/></w:rPr><w:t xml:space="preserve">property </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Menlo" w:hAnsi="Menlo" w:eastAsia="Menlo" w:cs="Menlo"/><w:b/><w:color w:val="0000FF"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:shd w:val="clear" w:fill="EFEFEF"/></w:rPr><w:t>name</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Menlo" w:hAnsi="Menlo" w:eastAsia="Menlo" w:cs="Menlo"/><w:b/><w:color w:val="008000"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:shd w:val="clear" w:fill="EFEFEF"/></w:rPr><w:t xml:space="preserve">="log.path" </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Menlo" w:hAnsi="Menlo" w:eastAsia="Menlo" w:cs="Menlo"/><w:b/><w:color w:val="0000FF"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:shd w:val="clear" w:fill="EFEFEF"/></w:rPr><w:t>value</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Menlo" w:hAnsi="Menlo" w:eastAsia="Menlo" w:cs="Menlo"/><w:b/><w:color w:val="008000"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:shd w:val="clear" w:fill="EFEFEF"/></w:rPr><w:t xml:space="preserve">="/Users/jyking/project/logs" </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Menlo" w:hAnsi="Menlo" w:eastAsia="Menlo" w:cs="Menlo"/><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:shd w:val="clear" w:fill="EFEFEF"/></w:rPr><w:t>/&gt;</w:t></w:r></w:p><w:p><w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/></w:p>
'@

$wrapped = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + '<w:body>' + $xmlPayload + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$newRange.InsertXML($wrapped)

Write-Output "done"
